# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets,
# reflecting the refreshed data output for the gh-pages build.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 7015
    $ws.Range("F3").Value = 56
    $ws.Range("F5").Value = 94
    $ws.Range("F6").Value = 1081

    if ($sheetName -eq "展览") {
        $ws.Range("F8").Value = 11
    } else {
        $ws.Range("F9").Value = 11
    }
}
